# Updates cryptocurrency price (D) and 1h volume-change (E) figures
# on Sheet1, matching the latest scrape pulled in by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.464.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.106.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5230"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4492"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.59%  "
$ws.Range("E9").Value = "  +16.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08953"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.100.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.767"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.773"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001125"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06609"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.302"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.522.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.346"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.347.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.203"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.677"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.166"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.937"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02571"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06794"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.498"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2283"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6929"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.255"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.316"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6382"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.642"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.248"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("E50").Value = "  +5.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.98%  "
